$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'71.332.90"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'  +2.63%  "
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(3,4).Value = "'4.000.90"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'  +1.51%  "
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(4,5).Value = "'  +0.15%  "
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(5,4).Value = "'529.56"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "'  +4.79%  "
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(6,4).Value = "'149.63"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "'  +1.30%  "
$ws.Cells.Item(6,5).Style = "Normal"
$ws.Cells.Item(7,5).Value = "'  -0.35%  "
$ws.Cells.Item(7,5).Style = "Normal"
$ws.Cells.Item(8,4).Value = "'0.999"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "'  +0.10%  "
$ws.Cells.Item(8,5).Style = "Normal"
$ws.Cells.Item(9,4).Value = "'0.737"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "'  +0.34%  "
$ws.Cells.Item(9,5).Style = "Normal"
$ws.Cells.Item(10,5).Value = "'  +0.24%  "
$ws.Cells.Item(10,5).Style = "Normal"
$ws.Cells.Item(11,4).Value = "'0.0000343"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "'  -2.11%  "
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(12,4).Value = "'43.18"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "'  -0.78%  "
$ws.Cells.Item(12,5).Style = "Normal"
$ws.Cells.Item(13,4).Value = "'10.65"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "'  +1.40%  "
$ws.Cells.Item(13,5).Style = "Normal"
$ws.Cells.Item(14,4).Value = "'4.637.45"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "'  +1.43%  "
$ws.Cells.Item(14,5).Style = "Normal"
$ws.Cells.Item(15,4).Value = "'4.005.75"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "'  +1.65%  "
$ws.Cells.Item(15,5).Style = "Normal"
$ws.Cells.Item(16,4).Value = "'21.33"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "'  +6.63%  "
$ws.Cells.Item(16,5).Style = "Normal"
$ws.Cells.Item(17,4).Value = "'14.34"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = "'  +0.72%  "
$ws.Cells.Item(17,5).Style = "Normal"
$ws.Cells.Item(18,4).Value = "'1.23"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "'  +2.02%  "
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(19,5).Value = "'  -1.88%  "
$ws.Cells.Item(19,5).Style = "Normal"
$ws.Cells.Item(20,4).Value = "'71.362.75"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "'  +2.69%  "
$ws.Cells.Item(20,5).Style = "Normal"
$ws.Cells.Item(21,4).Value = "'441.96"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "'  +1.27%  "
$ws.Cells.Item(21,5).Style = "Normal"
$ws.Cells.Item(22,4).Value = "'3.54"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "'  +2.82%  "
$ws.Cells.Item(22,5).Style = "Normal"
$ws.Cells.Item(23,4).Value = "'92.02"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "'  +3.37%  "
$ws.Cells.Item(23,5).Style = "Normal"
$ws.Cells.Item(24,4).Value = "'12.38"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "'  +3.29%  "
$ws.Cells.Item(24,5).Style = "Normal"
$ws.Cells.Item(25,4).Value = "'14.30"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "'  -2.83%  "
$ws.Cells.Item(25,5).Style = "Normal"
$ws.Cells.Item(26,4).Value = "'4.10"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "'  +5.78%  "
$ws.Cells.Item(26,5).Style = "Normal"
$ws.Cells.Item(27,4).Value = "'10.89"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "'  -2.98%  "
$ws.Cells.Item(27,5).Style = "Normal"
$ws.Cells.Item(28,4).Value = "'36.93"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = "'  -0.63%  "
$ws.Cells.Item(28,5).Style = "Normal"
$ws.Cells.Item(29,4).Value = "'13.62"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = "'  +1.07%  "
$ws.Cells.Item(29,5).Style = "Normal"
$ws.Cells.Item(30,4).Value = "'686.77"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = "'  -2.82%  "
$ws.Cells.Item(30,5).Style = "Normal"
$ws.Cells.Item(31,5).Value = "'  +0.19%  "
$ws.Cells.Item(31,5).Style = "Normal"
$ws.Cells.Item(33,4).Value = "'6.83"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "'  +12.67%  "
$ws.Cells.Item(33,5).Style = "Normal"
$ws.Cells.Item(34,4).Value = "'68.44"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "'  +5.76%  "
$ws.Cells.Item(34,5).Style = "Normal"
$ws.Cells.Item(35,4).Value = "'0.0₃0901"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = "'  +1.32%  "
$ws.Cells.Item(35,5).Style = "Normal"
$ws.Cells.Item(36,4).Value = "'0.444"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "'  -1.69%  "
$ws.Cells.Item(36,5).Style = "Normal"
$ws.Cells.Item(37,4).Value = "'40.91"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = "'  -0.11%  "
$ws.Cells.Item(37,5).Style = "Normal"
$ws.Cells.Item(38,5).Value = "'  +15.33%  "
$ws.Cells.Item(38,5).Style = "Normal"
$ws.Cells.Item(39,5).Value = "'  -0.84%  "
$ws.Cells.Item(39,5).Style = "Normal"
$ws.Cells.Item(40,5).Value = "'  +0.04%  "
$ws.Cells.Item(40,5).Style = "Normal"
$ws.Cells.Item(41,5).Value = "'  -0.01%  "
$ws.Cells.Item(41,5).Style = "Normal"
$ws.Cells.Item(42,5).Value = "'  +0.14%  "
$ws.Cells.Item(42,5).Style = "Normal"
$ws.Cells.Item(43,4).Value = "'2.90"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "'  +0.27%  "
$ws.Cells.Item(43,5).Style = "Normal"
$ws.Cells.Item(44,5).Value = "'  +0.38%  "
$ws.Cells.Item(44,5).Style = "Normal"
$ws.Cells.Item(45,4).Value = "'3.31"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "'  +9.83%  "
$ws.Cells.Item(45,5).Style = "Normal"
$ws.Cells.Item(46,4).Value = "'3.54"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "'  +5.09%  "
$ws.Cells.Item(46,5).Style = "Normal"
$ws.Cells.Item(47,4).Value = "'0.145"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "'  +0.61%  "
$ws.Cells.Item(47,5).Style = "Normal"
$ws.Cells.Item(48,4).Value = "'0.000285"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "'  +19.63%  "
$ws.Cells.Item(48,5).Style = "Normal"
$ws.Cells.Item(49,4).Value = "'9.30"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "'  +5.09%  "
$ws.Cells.Item(49,5).Style = "Normal"
$ws.Cells.Item(50,2).Value = "LidoDAOToken"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(50,4).Value = "'3.41"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "'  +0.18%  "
$ws.Cells.Item(50,5).Style = "Normal"
$ws.Cells.Item(51,2).Value = "BabyDogeCoin"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(51,4).Value = "'0.0₆0351"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "'  +0.08%  "
$ws.Cells.Item(51,5).Style = "Normal"
